# Update raw-score tables on both sheets ("060" and "066").

$wb = $excel.ActiveWorkbook

# --- Sheet "060" (first sheet) ---
$ws1 = $wb.Worksheets.Item("060")

$ws1.Range("B2").Value = 90
$ws1.Range("C2").Value = 86
$ws1.Range("B3").Value = 92
$ws1.Range("C3").Value = 91
$ws1.Range("B4").Value = 102
$ws1.Range("C4").Value = 97
$ws1.Range("B5").Value = 104
$ws1.Range("C5").Value = 108
$ws1.Range("B6").Value = 112
$ws1.Range("C6").Value = 111

# --- Sheet "066" (second sheet) ---
$ws2 = $wb.Worksheets.Item("066")

$ws2.Range("B2").Value = 85
$ws2.Range("C2").Value = 91
$ws2.Range("B3").Value = 92
$ws2.Range("C3").Value = 92
$ws2.Range("B4").Value = 96
$ws2.Range("C4").Value = 99
$ws2.Range("B5").Value = 105
$ws2.Range("C5").Value = 107
$ws2.Range("B6").Value = 111
$ws2.Range("C6").Value = 113

# Update the selection on sheet "066" to C7 (single cell, active cell C7).
$ws2.Activate()
$ws2.Range("C7").Select()
